$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target data for rows 2-13 (columns A-T), reflecting the refreshed TPM
# values and the addition of "Resolving-Mac" as a sending cluster.
$data = @(
    @("FAPs", "Fgf17", "Fgfr1", "ECs", 3, 1, 1.430739666666667, 4.292218999999999, 0.7311024264480693, 0.7311024264480693, 3, 1, 7.955277333333332, 23.865832, 0.05015625076675284, 0.05015625076675283, 11.38193084013422, 102.437377561208, 0.03666935663711084, 0.03666935663711084),
    @("FAPs", "Fgf17", "Fgfr1", "FAPs", 3, 1, 1.430739666666667, 4.292218999999999, 0.7311024264480693, 0.7311024264480693, 3, 1, 82.48060333333333, 247.44181, 0.520021823355633, 0.520021823355633, 118.0082709195989, 1062.07443827639, 0.3801892168612526, 0.3801892168612526),
    @("FAPs", "Fgf17", "Fgfr1", "MuSCs", 3, 1, 1.430739666666667, 4.292218999999999, 0.7311024264480693, 0.7311024264480693, 3, 1, 66.90297433333333, 200.708923, 0.4218083439585467, 0.4218083439585465, 95.72073919668188, 861.4866527701369, 0.3083851037641353, 0.3083851037641352),
    @("FAPs", "Fgf17", "Fgfr1", "Resolving-Mac", 3, 1, 1.430739666666667, 4.292218999999999, 0.7311024264480693, 0.7311024264480693, 3, 1, 1.271033333333333, 3.8131, 0.008013581919067616, 0.008013581919067614, 1.818517807655555, 16.3666602689, 0.00585874918557071, 0.005858749185570709),
    @("MuSCs", "Fgf17", "Fgfr1", "ECs", 2, 0.6666666666666666, 0.4223926666666666, 1.267178, 0.2158410161600821, 0.2158410161600822, 3, 1, 7.955277333333332, 23.865832, 0.05015625076675284, 0.05015625076675283, 3.360250806899555, 30.24225726209599, 0.01082577613227583, 0.01082577613227583),
    @("MuSCs", "Fgf17", "Fgfr1", "FAPs", 2, 0.6666666666666666, 0.4223926666666666, 1.267178, 0.2158410161600821, 0.2158410161600822, 3, 1, 82.48060333333333, 247.44181, 0.520021823355633, 0.520021823355633, 34.83920199024222, 313.55281791218, 0.1122420387784986, 0.1122420387784986),
    @("MuSCs", "Fgf17", "Fgfr1", "MuSCs", 2, 0.6666666666666666, 0.4223926666666666, 1.267178, 0.2158410161600821, 0.2158410161600822, 3, 1, 66.90297433333333, 200.708923, 0.4218083439585467, 0.4218083439585465, 28.25932573658822, 254.333931629294, 0.09104354158481416, 0.09104354158481415),
    @("MuSCs", "Fgf17", "Fgfr1", "Resolving-Mac", 2, 0.6666666666666666, 0.4223926666666666, 1.267178, 0.2158410161600821, 0.2158410161600822, 3, 1, 1.271033333333333, 3.8131, 0.008013581919067616, 0.008013581919067614, 0.5368751590888889, 4.8318764318, 0.001729659664493615, 0.001729659664493615),
    @("Resolving-Mac", "Fgf17", "Fgfr1", "ECs", 1, 0.3333333333333333, 0.1038296666666667, 0.311489, 0.05305655739184852, 0.05305655739184854, 3, 1, 7.955277333333332, 23.865832, 0.05015625076675284, 0.05015625076675283, 0.8259937937608888, 7.433944143848, 0.002661117997366169, 0.002661117997366169),
    @("Resolving-Mac", "Fgf17", "Fgfr1", "FAPs", 1, 0.3333333333333333, 0.1038296666666667, 0.311489, 0.05305655739184852, 0.05305655739184854, 3, 1, 82.48060333333333, 247.44181, 0.520021823355633, 0.520021823355633, 8.563933550565556, 77.07540195509, 0.02759056771588185, 0.02759056771588186),
    @("Resolving-Mac", "Fgf17", "Fgfr1", "MuSCs", 1, 0.3333333333333333, 0.1038296666666667, 0.311489, 0.05305655739184852, 0.05305655739184854, 3, 1, 66.90297433333333, 200.708923, 0.4218083439585467, 0.4218083439585465, 6.946513524038555, 62.51862171634701, 0.02237969860959721, 0.02237969860959721),
    @("Resolving-Mac", "Fgf17", "Fgfr1", "Resolving-Mac", 1, 0.3333333333333333, 0.1038296666666667, 0.311489, 0.05305655739184852, 0.05305655739184854, 3, 1, 1.271033333333333, 3.8131, 0.008013581919067616, 0.008013581919067614, 0.1319709673222222, 1.1877387059, 0.0004251730690032906, 0.0004251730690032906)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
    }
}
